$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Results" cell for the gdqmc_icf_pin2.0 (2/18/2018) block ---
$ws.Range("D76").Value = "calculated but don't know why not saved"

# --- Append new log entries (2018/3/26, cpqmc_bcs1.2 psudo-BCS runs) ---

# Row 118-119: cpqmc_bcs1.2DET u4
$ws.Range("A118").Value = "'3/26/2018"
$ws.Range("A118").ClearFormats()
$ws.Range("B118").Value = "cpqmc_bcs1.2DET u4"
$ws.Range("C118").Value = "~"
$ws.Range("E118").Value = "Bore"
$ws.Range("E119").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2DET_u4"

# Row 120-121: cpqmc_bcs1.2DET u8
$ws.Range("B120").Value = "cpqmc_bcs1.2DET u8"
$ws.Range("C120").Value = "~"
$ws.Range("E120").Value = "Bore"
$ws.Range("E121").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2DET_u8"

# Row 122-123: cpqmc_bcs1.2DET u12
$ws.Range("B122").Value = "cpqmc_bcs1.2DET u12"
$ws.Range("C122").Value = "~"
$ws.Range("E122").Value = "Bore"
$ws.Range("E123").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2DET_u12"

# Row 124-125: cpqmc_bcs1.2x u4
$ws.Range("B124").Value = "cpqmc_bcs1.2x u4"
$ws.Range("C124").Value = "~"
$ws.Range("E124").Value = "Bore"
$ws.Range("E125").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2x_u4"

# Row 126-127: cpqmc_bcs1.2x u8
$ws.Range("B126").Value = "cpqmc_bcs1.2x u8"
$ws.Range("C126").Value = "~"
$ws.Range("E126").Value = "Bore"
$ws.Range("E127").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2x_u8"

# Row 128-129: cpqmc_bcs1.2x u12
$ws.Range("B128").Value = "cpqmc_bcs1.2x u12"
$ws.Range("C128").Value = "~"
$ws.Range("E128").Value = "Bore"
$ws.Range("E129").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2x_u12"

# Row 130-131: cpqmc_bcs1.2xx u4
$ws.Range("B130").Value = "cpqmc_bcs1.2xx u4"
$ws.Range("C130").Value = "~"
$ws.Range("E130").Value = "Bore"
$ws.Range("E131").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2xx_u4"

# Row 132-133: cpqmc_bcs1.2xx u8
$ws.Range("B132").Value = "cpqmc_bcs1.2xx u8"
$ws.Range("C132").Value = "~"
$ws.Range("E132").Value = "Bore"
$ws.Range("E133").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2xx_u8"

# Row 134-135: cpqmc_bcs1.2xx u12
$ws.Range("B134").Value = "cpqmc_bcs1.2xx u12"
$ws.Range("C134").Value = "~"
$ws.Range("E134").Value = "Bore"
$ws.Range("E135").Value = "/sciclone/pscr/zxiao01/run/cpqmc_bcs1.2xx_u12"
